# Actualización automática 2025-06-25 13:20:09
#
# The underlying source number for "MOREIRA MOREIRA PATRICIO IGNACIO"
# increased by 622.08 (2851.2 -> 3473.28). This value is duplicated (as a
# cached/static snapshot, not a live formula) across the three report
# sheets, so every dependent cell has to be refreshed to match.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": raw detail row for the client ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("C17").Value = 3473.28

# --- Sheet "VENTA MENSUAL": monthly detail + column total ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F17").Value = 9485.16
$wsMensual.Range("F29").Value = 13815.49

# --- Sheet "CUMPLIMIENTO MENSUAL": per-group row + grand total ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D2").Value = 3473.28
$wsCumplimiento.Range("E2").Value = -3128.995395370514
$wsCumplimiento.Range("F2").Value = 10.08839766081871

$wsCumplimiento.Range("D19").Value = 20833.42
$wsCumplimiento.Range("E19").Value = 2666.580930050387
$wsCumplimiento.Range("F19").Value = 0.8865284755525042
